# Refresh the crypto price table (columns D = Price, E = Volume(1h)).
# Two rows (18/19 and 46/47) are also reordered - the coin that used to be
# second in the pair now comes first, each carrying its own refreshed
# price/volume.
#
# Price/Volume are stored as plain text (e.g. "26.519.76", "1.001",
# "0.000007950") rather than numbers, so for any value that Excel would
# otherwise auto-parse as a number we force the cell to Text format first
# and then assign the literal string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" '26.519.76'
$ws.Range("E2").Value = '  +0.12%  '

# Row 3 - Ethereum
Set-TextValue "D3" '1.843.11'
$ws.Range("E3").Value = '  -0.08%  '

# Row 4 - TetherUSD
Set-TextValue "D4" '1.001'
$ws.Range("E4").Value = '  +0.06%  '

# Row 5 - BNB
Set-TextValue "D5" '262.03'
$ws.Range("E5").Value = '  -0.33%  '

# Row 6 - USDC
$ws.Range("E6").Value = '  +0.06%  '

# Row 7 - XRP
Set-TextValue "D7" '0.5326'
$ws.Range("E7").Value = '  +2.43%  '

# Row 8 - Cardano
Set-TextValue "D8" '0.3060'
$ws.Range("E8").Value = '  -4.66%  '

# Row 9 - Dogecoin
Set-TextValue "D9" '0.06891'
$ws.Range("E9").Value = '  +1.68%  '

# Row 10 - Solana
Set-TextValue "D10" '18.26'
$ws.Range("E10").Value = '  -1.87%  '

# Row 11 - TRON
Set-TextValue "D11" '0.07793'
$ws.Range("E11").Value = '  +0.33%  '

# Row 12 - Polygon
Set-TextValue "D12" '0.7505'
$ws.Range("E12").Value = '  -2.58%  '

# Row 13 - WrappedEther
Set-TextValue "D13" '1.846.91'
$ws.Range("E13").Value = '  -0.06%  '

# Row 14 - Litecoin
Set-TextValue "D14" '89.74'
$ws.Range("E14").Value = '  +1.74%  '

# Row 15 - Polkadot
Set-TextValue "D15" '5.001'
$ws.Range("E15").Value = '  -0.10%  '

# Row 17 - Avalanche
$ws.Range("E17").Value = '  +0.65%  '

# Rows 18/19 - Dai <-> ShibaInu swap with refreshed price/volume
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D18" '0.000007950'
$ws.Range("E18").Value = '  +0.09%  '

$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D19" '1.001'
$ws.Range("E19").Value = '  +0.04%  '

# Row 20 - WrappedBTC
Set-TextValue "D20" '26.540.75'
$ws.Range("E20").Value = '  +0.05%  '

# Row 21 - WrappedliquidstakedEther2.0
Set-TextValue "D21" '2.078.51'
$ws.Range("E21").Value = '  -0.48%  '

# Row 22 - Uniswap
Set-TextValue "D22" '4.622'
$ws.Range("E22").Value = '  +0.33%  '

# Row 23 - Chainlink
Set-TextValue "D23" '5.990'
$ws.Range("E23").Value = '  +0.25%  '

# Row 24 - Cosmos
Set-TextValue "D24" '9.329'
$ws.Range("E24").Value = '  -1.01%  '

# Row 25 - Monero
Set-TextValue "D25" '143.03'
$ws.Range("E25").Value = '  -0.25%  '

# Row 26 - LidoDAOToken
Set-TextValue "D26" '2.199'
$ws.Range("E26").Value = '  +1.84%  '

# Row 27 - Toncoin
Set-TextValue "D27" '1.694'
$ws.Range("E27").Value = '  +1.13%  '

# Row 28 - EthereumClassic
Set-TextValue "D28" '16.98'
$ws.Range("E28").Value = '  +0.11%  '

# Row 29 - BitcoinCash
Set-TextValue "D29" '110.76'
$ws.Range("E29").Value = '  -0.50%  '

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" '4.269'
$ws.Range("E30").Value = '  +2.73%  '

# Row 31 - Stellar
$ws.Range("E31").Value = '  +1.09%  '

# Row 32 - Filecoin
Set-TextValue "D32" '4.074'
$ws.Range("E32").Value = '  -0.45%  '

# Row 33 - Hedera
Set-TextValue "D33" '0.04810'
$ws.Range("E33").Value = '  +0.00%  '

# Row 34 - HuobiToken
Set-TextValue "D34" '2.931'
$ws.Range("E34").Value = '  +2.43%  '

# Row 35 - ImmutableX
Set-TextValue "D35" '0.7281'
$ws.Range("E35").Value = '  +1.62%  '

# Row 36 - ARBITRUM
Set-TextValue "D36" '1.135'
$ws.Range("E36").Value = '  +0.99%  '

# Row 38 - RenderToken
Set-TextValue "D38" '2.302'
$ws.Range("E38").Value = '  +5.38%  '

# Row 39 - VeChain (price only; Volume(1h) unchanged)
Set-TextValue "D39" '0.01724'

# Row 40 - TheSandbox
Set-TextValue "D40" '0.4785'
$ws.Range("E40").Value = '  -0.74%  '

# Row 41 - TrustWalletToken
Set-TextValue "D41" '0.9107'
$ws.Range("E41").Value = '  +1.72%  '

# Row 42 - Quant
Set-TextValue "D42" '108.31'
$ws.Range("E42").Value = '  -3.27%  '

# Row 43 - FraxShare
Set-TextValue "D43" '5.880'
$ws.Range("E43").Value = '  -2.45%  '

# Row 44 - PaxDollar
$ws.Range("E44").Value = '  +0.07%  '

# Row 45 - Aptos
Set-TextValue "D45" '7.499'
$ws.Range("E45").Value = '  -1.10%  '

# Rows 46/47 - Decentraland <-> EnergySwap swap with refreshed price/volume
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D46" '9.097'
$ws.Range("E46").Value = '  +1.01%  '

$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue "D47" '0.4135'
$ws.Range("E47").Value = '  -0.66%  '

# Row 48 - Algorand
Set-TextValue "D48" '0.1242'
$ws.Range("E48").Value = '  +1.26%  '

# Row 49 - EOS
Set-TextValue "D49" '0.8998'
$ws.Range("E49").Value = '  +1.60%  '

# Row 50 - Elrond
Set-TextValue "D50" '34.83'
$ws.Range("E50").Value = '  -0.15%  '

# Row 51 - Cronos
Set-TextValue "D51" '0.05796'
$ws.Range("E51").Value = '  -1.74%  '
